# Update the ifo QoQ GVA error tables with recomputed statistics after
# adding a new forecast observation to each quarters-ahead (Q0..Q9) bucket.
# This reflects the ifo GDP component analysis preprocessing that folds in
# one additional matched error observation per horizon, shifting N (count)
# up by one and refreshing ME / MAE / MSE / RMSE / SE accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.02698403094954239, 0.5268716437348214, 0.4152143131774279, 0.644371254152005,  0.6664025289764616, 15)
    3  = @(0.2756751091950668,  0.5449100909190536, 0.4697505061201764, 0.6853834737722937, 0.6511854494255666, 14)
    4  = @(0.4115260714657983,  0.6969458320166075, 0.9427345523843811, 0.9709451850564897, 0.9153301309522154, 13)
    5  = @(0.6262380286368737,  0.6447007379214998, 0.8045669953357346, 0.8969765857232477, 0.6707333247272815, 12)
    6  = @(0.5245105352011469,  0.5475964834267181, 0.4278966426621968, 0.6541380914319215, 0.4099559430447661, 11)
    7  = @(0.5156605037261237,  0.5499064583805502, 0.4026603963166046, 0.6345552744376211, 0.3898071335267753, 10)
    8  = @(0.3926246975442547,  0.4457387173631203, 0.2339714725847577, 0.4837059774126817, 0.2996572782295139, 9)
    9  = @(0.5049883294650596,  0.5049883294650596, 0.284475404128805,  0.5333623572476829, 0.1880282677670355, 6)
    10 = @(0.451441112860743,   0.451441112860743,  0.2505375300224135, 0.5005372413940979, 0.2647785441877819, 3)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]  # B - ME
    $ws.Cells.Item($row, 3).Value = $values[1]  # C - MAE
    $ws.Cells.Item($row, 4).Value = $values[2]  # D - MSE
    $ws.Cells.Item($row, 5).Value = $values[3]  # E - RMSE
    $ws.Cells.Item($row, 6).Value = $values[4]  # F - SE
    $ws.Cells.Item($row, 7).Value = $values[5]  # G - N
}
